# Apply the "Alvearie FHIR IG" deploy update (v5.0.0 -> v6.0.0) to the
# StructureDefinition-line-of-business workbook.
#
# Sheet 1 "Metadata" is a Property/Value table describing the FHIR
# StructureDefinition. The old "Contact" metadata row (which just carried
# placeholder text) is replaced by a "Jurisdiction" row, the previously
# empty "Publisher" value is now filled in, the version/date are bumped,
# and one stale duplicate row is dropped so everything below shifts up by
# one row (21 rows -> 20 rows).
#
# Sheet 2 "Elements" is a details table for the StructureDefinition's
# elements; the only real content change there is the root Extension
# row's Short/Definition text, which used to be generic placeholder
# copy ("Extension" / "An Extension") and is now the real description.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item(1)   # "Metadata" sheet
$elements = $wb.Worksheets.Item(2)  # "Elements" sheet

# --- Metadata sheet -------------------------------------------------

# Drop the duplicate "Contact" row (old row 10); row 11 onward shifts up
# by one, which also removes the stray trailing "Context" row (old row
# 21) from the end of the table, taking the sheet from 21 to 20 rows.
$meta.Rows(10).Delete()

# Version 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date bump
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# New Jurisdiction row where the old duplicate "Contact" row used to be
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# --- Elements sheet --------------------------------------------------

# Root Extension row: Short / Definition go from generic placeholder
# text to the real (line-of-business-specific) description.
$elements.Cells.Item(2, 11).Value = "Line Of Business"
$elements.Cells.Item(2, 12).Value = "Code for the line of business"
